# Apply text-edits to the Missouri overview workbook factsheets.
#
# The "No. of 990 Filers w/ Gov Grants" counts were stored as real numbers
# on several tabs while every other metric column on those same tabs is
# stored as formatted text. This edit converts those count cells to text
# (using the classic leading-apostrophe trick so Excel treats the content
# literally instead of re-parsing it as a number), fills in proper
# percent/currency text for three counties that still had placeholder "0"
# values, and appends the missing statewide "Total" row to the County tab
# (mirroring the Total row already present on the other tabs).

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Leading apostrophe forces Excel to store the literal text instead of
    # reinterpreting look-like-a-number / look-like-a-percent content.
    $range.Value = "'" + $text
    # Clear the quote-prefix formatting flag Excel applies automatically
    # when the apostrophe trick is used, so no stray style is left behind.
    $range.Style = "Normal"
}

# ----------------------------------------------------------------------
# Overall sheet: A2 (filer count) -> text "1,959"
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overall")
Set-TextValue $ws.Range("A2") "1,959"

# ----------------------------------------------------------------------
# County sheet: column B (filer counts), rows 2-109 -> text
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("County")
$countyCounts = @("17","3","4","5","8","2","3","1","2","84","42","10","1","11","10","34","6","2","13","2","1","9","3","37","3","44","6","2","1","4","2","2","4","2","9","17","6","3","97","5","3","6","2","1","4","13","4","322","27","14","10","7","11","3","4","6","4","10","4","3","14","1","2","5","1","4","6","2","8","10","1","2","4","7","13","17","4","24","7","9","2","12","4","5","3","12","1","3","13","3","3","35","2","9","411","260","6","3","6","1","17","7","4","7","2","4","4","4")
for ($i = 0; $i -lt $countyCounts.Length; $i++) {
    $row = $i + 2
    Set-TextValue $ws.Cells.Item($row, 2) $countyCounts[$i]
}

# Rows 110-112 (Maries, McDonald, Ozark counties) had placeholder "0"
# values in every metric column; they become formatted percent/currency
# text like every other county row.
$zeroRows = @(110, 111, 112)
foreach ($r in $zeroRows) {
    Set-TextValue $ws.Cells.Item($r, 2) "0.00%"
    Set-TextValue $ws.Cells.Item($r, 3) "`$0"
    Set-TextValue $ws.Cells.Item($r, 4) "0.00%"
    Set-TextValue $ws.Cells.Item($r, 5) "0.00%"
    Set-TextValue $ws.Cells.Item($r, 6) "0.00%"
}

# New row 113: statewide Total row (matches the Total row already present
# on the other factsheet tabs).
Set-TextValue $ws.Range("A113") "Total"
Set-TextValue $ws.Range("B113") "1,959"
Set-TextValue $ws.Range("C113") "`$3,985,106,420"
Set-TextValue $ws.Range("D113") "7.82%"
Set-TextValue $ws.Range("E113") "-12.22%"
Set-TextValue $ws.Range("F113") "68.25%"

# ----------------------------------------------------------------------
# Congressional District sheet: column B (filer counts), rows 2-9 -> text,
# plus the existing Total row (10).
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Congressional District")
$cdCounts = @("526","172","148","183","341","228","176","185")
for ($i = 0; $i -lt $cdCounts.Length; $i++) {
    $row = $i + 2
    Set-TextValue $ws.Cells.Item($row, 2) $cdCounts[$i]
}
Set-TextValue $ws.Range("B10") "1,959"

# ----------------------------------------------------------------------
# Size sheet: column B (filer counts), rows 2-7 -> text, plus the existing
# Total row (8).
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Size")
$sizeCounts = @("535","511","294","132","394","93")
for ($i = 0; $i -lt $sizeCounts.Length; $i++) {
    $row = $i + 2
    Set-TextValue $ws.Cells.Item($row, 2) $sizeCounts[$i]
}
Set-TextValue $ws.Range("B8") "1,959"

# ----------------------------------------------------------------------
# Subsector sheet: column B (filer counts), rows 2-12 -> text, plus the
# existing Total row (13).
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Subsector")
$subsectorCounts = @("154","135","59","204","69","652","12","172","50","426","26")
for ($i = 0; $i -lt $subsectorCounts.Length; $i++) {
    $row = $i + 2
    Set-TextValue $ws.Cells.Item($row, 2) $subsectorCounts[$i]
}
Set-TextValue $ws.Range("B13") "1,959"
